# [TEST SCRAPE] updated files from azure vm
#
# 1. Clear the stray empty cell B6 on the "ODI Batting" sheet.
# 2. Add a new worksheet "ODI Batting Extra" at the end of the workbook
#    holding extra per-innings batting stats (4s/6s/percent of team runs/
#    man-of-the-match flag).

$wb = $excel.ActiveWorkbook

# --- 1. Clear B6 on "ODI Batting" ---------------------------------------
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$odiBatting.Range("B6").Value = ""

# --- 2. Add the "ODI Batting Extra" sheet at the end ---------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$extra = $wb.Worksheets.Add($null, $lastSheet)
$extra.Name = "ODI Batting Extra"

# Format-donor cells already present in the workbook:
#  - a header cell (bold / bordered / centered) to stamp on row 1
#  - a plain unstyled cell to stamp on ordinary cells, so that numeric-
#    looking strings ("4472", "3.55%", ...) are stored as real text instead
#    of being auto-coerced to numbers/percentages by Excel, without leaving
#    behind a stray "@" text-format style.
$headerDonor = $wb.Worksheets.Item("Player Info").Range("A1")
$plainDonor  = $wb.Worksheets.Item("Player Info").Range("A2")

function Set-TextCell($cell, [string]$text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $plainDonor.Copy()
    $cell.PasteSpecial(-4122) | Out-Null
}

# Header row ---------------------------------------------------------------
$headers = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($col = 1; $col -le $headers.Length; $col++) {
    $cell = $extra.Cells.Item(1, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $headers[$col - 1]
    $headerDonor.Copy()
    $cell.PasteSpecial(-4122) | Out-Null
}

# Data rows ------------------------------------------------------------
# MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
# ($null marks the blank cells in the trailing "did not bat" row.)
$rows = @(
    @("4472", 5, "1", "0", "3.55%", "NO"),
    @("4473", 5, "4", "0", "28.72%", "NO"),
    @("4476", 9, "0", "0", "0.91%", "NO"),
    @("4564", 4, "0", "0", "1.33%", "NO"),
    @("4565", $null, $null, $null, $null, "NO")
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowNum = $r + 2
    $row = $rows[$r]

    Set-TextCell $extra.Cells.Item($rowNum, 1) $row[0]

    if ($null -ne $row[1]) {
        $extra.Cells.Item($rowNum, 2).Value = $row[1]
    } else {
        Set-TextCell $extra.Cells.Item($rowNum, 2) ""
    }

    if ($null -ne $row[2]) {
        Set-TextCell $extra.Cells.Item($rowNum, 3) $row[2]
    } else {
        Set-TextCell $extra.Cells.Item($rowNum, 3) ""
    }

    if ($null -ne $row[3]) {
        Set-TextCell $extra.Cells.Item($rowNum, 4) $row[3]
    } else {
        Set-TextCell $extra.Cells.Item($rowNum, 4) ""
    }

    if ($null -ne $row[4]) {
        Set-TextCell $extra.Cells.Item($rowNum, 5) $row[4]
    } else {
        Set-TextCell $extra.Cells.Item($rowNum, 5) ""
    }

    Set-TextCell $extra.Cells.Item($rowNum, 6) $row[5]
}
